$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,6,8,9,10,11,12,14,16,17,18,21,22,23,24,25,27,28,29,31,32,33,34,38,39,40,41,42,43,45,46,47,49,52,54,55)

foreach ($r in $rows) {
    $ws.Range("B$r").Value = "No"
    $ws.Range("F$r").Value = ""
    $ws.Range("J$r").Value = 0
}
